$d = $word.ActiveDocument

# 1) Merge the three runs ("...e-" / "commerce" / ".") that together
#    spell "...e-commerce." into a single run and drop the spell-check
#    proofErr markers that wrapped "commerce". A Find/Replace whose
#    replacement text equals the found text makes Word rewrite the
#    matched range as a single run.
$d.Content.Find.Execute(
    "sistema e-commerce.", $true, $false, $false, $false, $false,
    $true, 1, $false, "sistema e-commerce.", 2
) | Out-Null

# 2) Locate the paragraph that ends section 2.5's list ("Hacer envíos...")
#    and append a new paragraph right after it with the closing point
#    about handling defective products.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Hacer envíos a toda la república mexicana*") {
        $target = $p
        break
    }
}

$target.Range.InsertParagraphAfter()

$newParaIndex = $target.Index + 1
$newPara = $d.Paragraphs.Item($newParaIndex)
$newPara.Range.Text = "Implementar un área dedicada al recibo de productos defectuosos de fabrica , que recopilen la información del producto , del usuario, que den seguimiento al reenvió del mismo, confirmen que el producto tiene un error de fabrica o problema y hacer el reembolso de la cantidad aceptada o entregar otro producto con el mismo valor."
